# The source feed re-stamped its "Förändrad" (last-changed) column on every
# row of the sheet: 45188 (2023-09-19) -> 45189 (2023-09-20).
# Update column C for every data row (row 2 through the last used row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45188) {
        $cell.Value = 45189
    }
}
